$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data.
# Values are written with a leading apostrophe to force text
# interpretation (the source cells are plain text, and many of the
# new values look numeric, e.g. "587.57" or "63.483.52", which Excel
# would otherwise silently convert to a Number). The style is reset
# back to Normal immediately after so no stray quote-prefix formatting
# is left behind on the cell.

$ws.Range('D2').Value = "'63.483.52"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -1.38%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.600.17"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -1.65%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.00%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'587.57"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -2.98%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'149.26"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -1.83%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.00%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.583"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -1.55%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.109"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -1.35%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  +1.55%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.386"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -1.27%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  -0.67%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'27.50"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -1.12%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'3.065.32"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -1.70%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'63.289.90"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -1.43%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Value = "'  +3.35%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'2.632.60"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -0.50%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'12.05"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -1.62%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'4.66"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -0.20%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'343.98"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -2.64%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'  -2.34%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'  +0.07%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'66.48"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.66%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'1.72"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -1.98%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  -2.40%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'1.65"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -4.41%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'564.13"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +2.98%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'8.19"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -0.57%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('B29').Value = "'Kaspa"
$ws.Range('B29').Style = 'Normal'
$ws.Range('C29').Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('C29').Style = 'Normal'
$ws.Range('D29').Value = "'0.161"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -3.54%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('B30').Value = "'Binance-PegBSC-USD"
$ws.Range('B30').Style = 'Normal'
$ws.Range('C30').Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range('C30').Style = 'Normal'
$ws.Range('D30').Value = "'0.999"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -0.03%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D32').Value = "'0.0₃0843"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -3.36%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'1.76"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -0.76%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'5.28"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -0.97%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'165.56"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -1.23%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'0.412"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -0.02%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  -0.05%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'19.40"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -1.25%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'1.92"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -5.36%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  -0.07%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'165.71"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -1.48%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'3.97"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +0.53%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'22.86"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +5.08%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.0579"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -1.59%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'2.11"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +1.83%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.630"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -0.06%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'0.0247"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -0.16%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'0.0957"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -1.06%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'19.05"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -1.83%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  +13.57%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'0.179"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -3.93%  "
$ws.Range('E51').Style = 'Normal'
